$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header columns to the new catalog vocabulary.
$ws.Range("B1").Value = "t1"
$ws.Range("C1").Value = "t2"
$ws.Range("D1").Value = "effect"

# Record which effect measure is being reported (new annotation cell).
$ws.Range("I2").Value = "effect=md"
$ws.Range("I2").Font.Bold = $true

# Leave the selection where the editor finished working.
$ws.Range("G7").Select() | Out-Null
